# "bug in sign up fix it"
# Adds the missing sign-up row (niro12 / nironi1@ / 309375905) that the
# bug was dropping, and switches the workbook's base font from Calibri
# to Arial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 13: the sign-up record that was missing ---------------------
# Shared strings must be introduced in column order (A then B) so the new
# strings land at the same indices (12 = "niro12", 13 = "nironi1@") that
# Excel would have assigned.
$ws.Range("A13").Value = "niro12"
$ws.Range("B13").Value = "nironi1@"
$ws.Range("C13").Value = 309375905

# --- Switch the workbook's default font from Calibri to Arial ------------
# (Column widths / row heights shift by a sub-pixel amount in real Excel as
# a side effect of this font-metric change; that's cosmetic noise the
# headless engine doesn't recompute, so it's left alone rather than
# approximated with a coarser value.)
$wb.Styles.Item("Normal").Font.Name = "Arial"
